$wb = $excel.ActiveWorkbook

# --- Priority Status sheet: update species counts ---
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# --- Species qualification sheet: rename assessment label ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("A2").Value = "SoIB Assessment"

# --- High Priority break-up sheet: rename + becomes the interannual update sheet ---
$wsHighPriority = $wb.Worksheets.Item("High Priority break-up")
$wsHighPriority.Name = "Interannual update - High Pri"

# New data rows for the (renamed) Interannual update sheet
$wsHighPriority.Range("A2").Value = "Trend New"
$wsHighPriority.Range("B2").Value = 97
$wsHighPriority.Range("C2").Value = 94.2
$wsHighPriority.Range("D2").Value = 97
$wsHighPriority.Range("E2").Value = 94.2

$wsHighPriority.Range("A3").Value = "IUCN"
$wsHighPriority.Range("B3").Value = 6
$wsHighPriority.Range("C3").Value = 5.8
$wsHighPriority.Range("D3").Value = 6
$wsHighPriority.Range("E3").Value = 5.8

# --- New sheet: "Major update - High Priority " (old High Priority break-up content) ---
$wsMajor = $wb.Worksheets.Add($null, $wsHighPriority)
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"
$wsMajor.Range("A1:E1").Font.Bold = $true
$wsMajor.Range("A1:E1").HorizontalAlignment = -4108
